$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1:G41").AutoFilter(5, @("Phát triển phần mềm", ""), 7)

$n = $ws.Names.Add("_xlnm._FilterDatabase", "=Sheet1!`$A`$1:`$G`$41")
$n.Visible = $false

$ws.Range("G24").Select()
